# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# figures for this hours snapshot, as produced by the scheduled GitHub
# Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.422.29"
$ws.Range("D3").Value = "1.848.11"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.78"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6279"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07678"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2919"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.73"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.843.68"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.028"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001075"
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6799"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.46"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.172"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "29.448.48"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.01"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.07"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.394"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.68"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.345"
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05693"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.115"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.023"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.841"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.161"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7082"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.583"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.775"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").Value = "1.228.80"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.541"
$ws.Range("E40").Value = "  +4.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9109"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "2.002.17"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.74"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.03"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.137"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4013"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.984"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.669"
$ws.Range("E51").Value = "  +0.07%  "
